$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 196
$ws.Range("C7").Value = 196

$ws.Range("E6").Select()
